$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sprint1 sheet ("Sprint1") - re-colour column-A story-id cells red (matches
# existing direct-formatting style already used elsewhere in the workbook).
# A2/A3 additionally keep their wrap-text formatting (-> new cellXf with
# wrapText), A4:A11 simply pick up the existing "red text" cellXf (index 17
# in the original file).
# ---------------------------------------------------------------------------
$sprint1 = $wb.Worksheets.Item("Sprint1")
$sprint1.Range("A2:A3").Font.Color = 255
$sprint1.Range("A2:A3").WrapText = $true
$sprint1.Range("A4:A11").Font.Color = 255

# ---------------------------------------------------------------------------
# Sprint2 sheet - same red-text treatment for the story-id cells that need
# it (A2:A6 and A10); the others (A7:A9, A11) are left on their old style.
# ---------------------------------------------------------------------------
$sprint2 = $wb.Worksheets.Item("Sprint2")
$sprint2.Range("A2:A6").Font.Color = 255
$sprint2.Range("A10").Font.Color = 255

# ---------------------------------------------------------------------------
# Sprint3 sheet - add the 4 new backlog rows (US10, US19, US20, US23) that
# were completed this sprint. Formatting is pulled from a sister row on the
# Sprint2 sheet (same column layout/styles: s="30" story id, s="32" wrapped
# description, default style for the PT/SA column) via copy/paste-special of
# formats only, then the literal values are written in.
# ---------------------------------------------------------------------------
$sprint3 = $wb.Worksheets.Item("Sprint3")

$sprint2.Range("A7:C7").Copy()
$sprint3.Range("A2").PasteSpecial(-4122)
$sprint3.Range("A2").Value = "US10"
$sprint3.Range("B2").Value = "Marriage after 14"
$sprint3.Range("C2").Value = "PT"

$sprint2.Range("A8:C8").Copy()
$sprint3.Range("A3").PasteSpecial(-4122)
$sprint3.Range("A3").Value = "US19"
$sprint3.Range("B3").Value = "First cousins should not marry"
$sprint3.Range("C3").Value = "SA"

$sprint2.Range("A9:C9").Copy()
$sprint3.Range("A4").PasteSpecial(-4122)
$sprint3.Range("A4").Value = "US20"
$sprint3.Range("B4").Value = "Aunts and uncles"
$sprint3.Range("C4").Value = "SA"

$sprint2.Range("A11:C11").Copy()
$sprint3.Range("A5").PasteSpecial(-4122)
$sprint3.Range("A5").Value = "US23"
$sprint3.Range("B5").Value = "Unique name and birth date"
$sprint3.Range("C5").Value = "SA"

$sprint3.Rows.Item(2).RowHeight = 26
$sprint3.Rows.Item(3).RowHeight = 52
$sprint3.Rows.Item(4).RowHeight = 26
$sprint3.Rows.Item(5).RowHeight = 39

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Window / selection state: Sprint1 -> A9, Sprint2 -> A11:C11, Sprint3 ->
# A5:C5 (Sprint3 selected last so it ends up the active sheet/tab, matching
# the bumped activeTab index in workbook.xml).
# ---------------------------------------------------------------------------
$sprint1.Range("A9").Select()
$sprint2.Range("A11:C11").Select()
$sprint3.Range("A5:C5").Select()
